$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update progress percentages in column G (rows 6-19)
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 1
$ws.Range("G8").Value = 1
$ws.Range("G9").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("G11").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
$ws.Range("G14").Value = 0.3
$ws.Range("G15").Value = 0.8
$ws.Range("G16").Value = 1
$ws.Range("G17").Value = 1
$ws.Range("G18").Value = 1
$ws.Range("G19").Value = 1

# Update the view: scroll back to A1 (remove topLeftCell="D1") and select I14
$ws.Range("A1").Select() | Out-Null
$ws.Range("I14").Select() | Out-Null
